$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.887.14'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '1.647.31'
$ws.Range("E3").Value = '  +1.99%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''309.05'
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '''0.3889'
$ws.Range("E7").Value = '  -0.87%  '
$ws.Range("D8").Value = '''0.3824'
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = '''51.22'
$ws.Range("E9").Value = '  +4.42%  '
$ws.Range("D10").Value = '''1.347'
$ws.Range("E10").Value = '  -1.53%  '
$ws.Range("D11").Value = '''1.001'
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").Value = '''0.08426'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").Value = '''23.82'
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("D14").Value = '''7.090'
$ws.Range("E14").Value = '  +0.99%  '
$ws.Range("D15").Value = '''7.784'
$ws.Range("E15").Value = '  +3.33%  '
$ws.Range("D16").Value = '''0.00001308'
$ws.Range("E16").Value = '  +2.80%  '
$ws.Range("D17").Value = '1.648.48'
$ws.Range("E17").Value = '  +2.76%  '
$ws.Range("D18").Value = '''94.55'
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").Value = '''0.06977'
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("D20").Value = '''19.66'
$ws.Range("E20").Value = '  -1.67%  '
$ws.Range("D21").Value = '''6.856'
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("D22").Value = '''1.000'
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("D24").Value = '23.890.89'
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("D25").Value = '''2.469'
$ws.Range("E25").Value = '  +1.00%  '
$ws.Range("D26").Value = '''2.980'
$ws.Range("E26").Value = '  +4.28%  '
$ws.Range("D27").Value = '''21.99'
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").Value = '''152.38'
$ws.Range("E28").Value = '  -2.77%  '
$ws.Range("D29").Value = '''5.427'
$ws.Range("E29").Value = '  +3.84%  '
$ws.Range("D30").Value = '''138.77'
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("D31").Value = '''7.713'
$ws.Range("E31").Value = '  -1.79%  '
$ws.Range("D32").Value = '''2.485'
$ws.Range("E32").Value = '  -0.15%  '
$ws.Range("D33").Value = '1.830.20'
$ws.Range("E33").Value = '  +2.32%  '
$ws.Range("D34").Value = '''1.024'
$ws.Range("E34").Value = '  +5.07%  '
$ws.Range("D35").Value = '''0.08023'
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("D36").Value = '''0.02953'
$ws.Range("E36").Value = '  +2.80%  '
$ws.Range("D37").Value = '''6.692'
$ws.Range("E37").Value = '  +2.11%  '
$ws.Range("D38").Value = '''10.79'
$ws.Range("E38").Value = '  +5.51%  '
$ws.Range("D39").Value = '''0.2677'
$ws.Range("E39").Value = '  +0.82%  '
$ws.Range("D40").Value = '''0.09087'
$ws.Range("D41").Value = '''0.7513'
$ws.Range("E41").Value = '  +0.58%  '
$ws.Range("E42").Value = '  -0.40%  '
$ws.Range("D43").Value = '''1.415'
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("D44").Value = '''16.18'
$ws.Range("E44").Value = '  +1.85%  '
$ws.Range("D45").Value = '''0.6897'
$ws.Range("E45").Value = '  +0.49%  '
$ws.Range("D46").Value = '''2.438'
$ws.Range("E46").Value = '  -0.92%  '
$ws.Range("D47").Value = '''4.072'
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").Value = '''0.08279'
$ws.Range("E49").Value = '  +0.66%  '
$ws.Range("D50").Value = '''134.19'
$ws.Range("E50").Value = '  +0.72%  '
$ws.Range("D51").Value = '''1.215'
$ws.Range("E51").Value = '  +0.91%  '
